# Weekly fruit/vegetable data refresh: the price records in rows 2-13 are
# re-shuffled onto different rows (dates + volumes + prices move together as
# whole records), while the row 1 header and the non-moving descriptive
# columns (A,B,C,E,F,G,H,I,R) stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for the columns that actually change: D (Fecha),
# J (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio
# ponderado), N (Unidad de comercializacion), O (Origen), P (Precio $/Kg),
# Q (Kg o Unidades).
$rows = @{
    2  = @{ D = 44277; J = 25; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 167; Q = 60 }
    3  = @{ D = 44179; J = 15; K = 7000;  L = 7000;  M = 7000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 117; Q = 60 }
    4  = @{ D = 44284; J = 35; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 167; Q = 60 }
    5  = @{ D = 44585; J = 30; K = 11000; L = 11000; M = 11000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 183; Q = 60 }
    6  = @{ D = 44312; J = 30; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 167; Q = 60 }
    7  = @{ D = 44315; J = 25; K = 10000; L = 10000; M = 10000; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 167; Q = 60 }
    8  = @{ D = 44333; J = 25; K = 10000; L = 11000; M = 10400; N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 173; Q = 60 }
    9  = @{ D = 44200; J = 10; K = 9000;  L = 9000;  M = 9000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 150; Q = 60 }
    10 = @{ D = 44243; J = 80; K = 10000; L = 11000; M = 10375; N = "`$/caja 60 unidades"; O = "Provincia de Quillota"; P = 173; Q = 60 }
    11 = @{ D = 44405; J = 45; K = 9000;  L = 9000;  M = 9000;  N = "`$/caja 50 unidades"; O = "Provincia de Quillota"; P = 180; Q = 50 }
    12 = @{ D = 44186; J = 15; K = 7000;  L = 7000;  M = 7000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 117; Q = 60 }
    13 = @{ D = 44291; J = 20; K = 9000;  L = 9000;  M = 9000;  N = "`$/caja 60 unidades"; O = "Provincia de Limarí"; P = 150; Q = 60 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
}
